# Auto-generated Excel COM-interop script to apply diff changes
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws4 = $wb.Worksheets.Item(4)

$ws1.Cells.Item(4, 6).Value = 718
$ws1.Cells.Item(7, 6).Value = 2789
$ws1.Cells.Item(8, 6).Value = 1672
$ws1.Cells.Item(9, 6).Value = 1782
$ws1.Cells.Item(10, 6).Value = 317
$ws1.Cells.Item(12, 6).Value = 722
$ws1.Cells.Item(13, 6).Value = 884
$ws1.Cells.Item(14, 6).Value = 162
$ws1.Cells.Item(15, 6).Value = 363
$ws1.Cells.Item(16, 6).Value = 1111
$ws1.Cells.Item(18, 6).Value = 48
$ws1.Cells.Item(18, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/HHHVcvcC1709015213282.png"
$ws1.Cells.Item(20, 6).Value = 6472
$ws1.Cells.Item(22, 6).Value = 1441
$ws1.Cells.Item(23, 6).Value = 144
$ws1.Cells.Item(25, 6).Value = 157
$ws1.Cells.Item(26, 6).Value = 302
$ws1.Cells.Item(27, 6).Value = 259
$ws1.Cells.Item(28, 6).Value = 62
$ws1.Cells.Item(29, 6).Value = 1097
$ws1.Cells.Item(30, 6).Value = 893
$ws1.Cells.Item(32, 6).Value = 86
$ws1.Cells.Item(34, 6).Value = 467
$ws1.Cells.Item(35, 6).Value = 1361
$ws1.Cells.Item(41, 6).Value = 184
$ws1.Cells.Item(42, 6).Value = 152
$ws4.Cells.Item(4, 6).Value = 718
$ws4.Cells.Item(10, 6).Value = 2789
$ws4.Cells.Item(11, 6).Value = 1672
$ws4.Cells.Item(12, 6).Value = 1782
$ws4.Cells.Item(13, 6).Value = 317
$ws4.Cells.Item(15, 6).Value = 722
$ws4.Cells.Item(17, 6).Value = 884
$ws4.Cells.Item(18, 6).Value = 162
$ws4.Cells.Item(19, 6).Value = 363
$ws4.Cells.Item(20, 6).Value = 1111
$ws4.Cells.Item(21, 6).Value = 48
$ws4.Cells.Item(21, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/HHHVcvcC1709015213282.png"
$ws4.Cells.Item(23, 6).Value = 6472
$ws4.Cells.Item(25, 6).Value = 1441
$ws4.Cells.Item(27, 6).Value = 144
$ws4.Cells.Item(29, 6).Value = 157
$ws4.Cells.Item(30, 6).Value = 302
$ws4.Cells.Item(31, 6).Value = 259
$ws4.Cells.Item(32, 6).Value = 62
$ws4.Cells.Item(33, 6).Value = 1097
$ws4.Cells.Item(34, 6).Value = 893
$ws4.Cells.Item(36, 6).Value = 86
$ws4.Cells.Item(38, 6).Value = 467
$ws4.Cells.Item(39, 6).Value = 1361
$ws4.Cells.Item(45, 6).Value = 184
$ws4.Cells.Item(49, 6).Value = 152
